$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ballots")

# Row 17: Jay Dunn ballot
$ws.Range("A17").Value = "Jay Dunn"
$ws.Range("E17").Value = "x"
$ws.Range("I17").Value = "x"
$ws.Range("J17").Value = "x"
$ws.Range("K17").Value = "x"
$ws.Range("M17").Value = "x"
$ws.Range("O17").Value = "x"
$ws.Range("Q17").Value = "x"
$ws.Range("T17").Value = "x"
$ws.Range("AK17").Value = 8
$ws.Range("AL17").Value = "The Trentonian"
$ws.Range("AM16").Copy($ws.Range("AM17"))
$ws.Range("AM17").Value = 43439

# Row 18: Bob Herzog ballot
$ws.Range("A18").Value = "Bob Herzog"
$ws.Range("C18").Value = "x"
$ws.Range("D18").Value = "x"
$ws.Range("E18").Value = "x"
$ws.Range("F18").Value = "x"
$ws.Range("H18").Value = "x"
$ws.Range("I18").Value = "x"
$ws.Range("J18").Value = "x"
$ws.Range("K18").Value = "x"
$ws.Range("O18").Value = "x"
$ws.Range("V18").Value = "x"
$ws.Range("AK18").Value = 10
$ws.Range("AL18").Value = "Twitter"
$ws.Range("AM16").Copy($ws.Range("AM18"))
$ws.Range("AM18").Value = 43439

# Update the view: scroll so column AB is the left-most visible column,
# and select H17 as the active cell.
$ws.Range("H17").Select()
$excel.ActiveWindow.ScrollColumn = $ws.Range("AB1").Column
